# 0830 transparency sub coding
# Update codebook about coding subelements of transparency index.
#
# This script:
#  1. Relocates the existing "Fashion Transparency Index" block of rows
#     (old rows 43-45 and 47-51) down by 23 rows, to rows 66-68 and 70-74,
#     leaving rows 52-62 untouched in their original location.
#  2. Adds the new "transparency sub coding" rows 77-98 describing the
#     2018/2019 buyer transparency sub-indices (policy, governance,
#     traceability, know/audit-show-fix, etc).
#  3. Updates the active window selection to match the author's final
#     cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: move the FTI block (A43:C51, which includes one blank row at
# 46) down to A66:C74, preserving styles/values, then clear the old
# location and any left-over phantom cells created by the block copy.
# ---------------------------------------------------------------------
$src = $ws.Range("A43:C51")
$dst = $ws.Range("A66")
$src.Copy($dst)

$ws.Range("A43:C51").Clear()

# the source rows did not have data in column C for rows 44/45 (-> 67/68),
# and row 46 (-> 69) was fully blank; remove the phantom cells the block
# copy created there so the row disappears entirely, matching the source
# layout exactly.
$ws.Range("C67").Clear()
$ws.Range("C68").Clear()
$ws.Range("A69:C69").Clear()

# ---------------------------------------------------------------------
# Step 2: add the new rows describing the transparency sub-indices.
# ---------------------------------------------------------------------
function Set-RedLabel($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Color = $ws.Range("A57").Font.Color()
}

Set-RedLabel "A77" "buyer1pol18"
Set-RedLabel "B77" 'the first buyer''s transparency index regarding "policy and commitment" in 2018.'
$ws.Range("C77").Value = 'Values: 0=unrated or missing buyer names; 1 for index score between 0-10%, 2 for 11-20, 3 for 21-30, 4 for 31-40, 5 for 41-50, 6 for 51-60, 7 for 61-70, 8 for 71-80, 9 for 81-90, 10 for 91-100. '

Set-RedLabel "A78" "buyer1gov18"
Set-RedLabel "B78" 'the first buyer''s transparency index regarding "governance" in 2018.'
$ws.Range("C78").Value = "values: same as above "

Set-RedLabel "A79" "buyer1trace18"
Set-RedLabel "B79" 'the first buyer''s transparency index regarding "traceability" in 2018.'
$ws.Range("C79").Value = "values: same as above "

Set-RedLabel "A80" "buyer2trace18"
Set-RedLabel "B80" 'the second buyer''s transparency index regarding "traceability" in 2018.'
$ws.Range("C80").Value = "values: same as above "

Set-RedLabel "A81" "buyer3trace18"
Set-RedLabel "B81" 'the third buyer''s transparency index regarding "traceability" in 2018.'
$ws.Range("C81").Value = "values: same as above "

Set-RedLabel "A82" "buyertracemx18"
Set-RedLabel "B82" "the maximum traceability index among the three buyers in 2018"
$ws.Range("C82").Value = "values: same as above "

Set-RedLabel "A83" "buyer1know18"
Set-RedLabel "B83" 'the first buyer''s transparency index regarding "know/audit, show, and fix" in 2018'
$ws.Range("C83").Value = "values: same as above "

Set-RedLabel "A84" "buyer2know18"
Set-RedLabel "B84" 'the second buyer''s transparency index regarding "know/audit, show, and fix" in 2018'
$ws.Range("C84").Value = "values: same as above "

Set-RedLabel "A85" "buyer3know18"
Set-RedLabel "B85" 'the third buyer''s transparency index regarding "know/audit, show, and fix" in 2018'
$ws.Range("C85").Value = "values: same as above "

Set-RedLabel "A86" "buyerknowmx18"
Set-RedLabel "B86" 'the biggest "know/audit, show, and fix" index among three buyers in 2018'
$ws.Range("C86").Value = "values: same as above "

# blank divider row, still red-styled like the block above/below it
Set-RedLabel "A87" ""
Set-RedLabel "B87" ""

Set-RedLabel "A88" "the same set of variables and coding for the transparency index report in 2019 (which included 50 more brands). "

Set-RedLabel "A89" "buyer1pol19"
Set-RedLabel "B89" "definition is the same as 2018"

Set-RedLabel "A90" "buyer1gov19"
Set-RedLabel "A91" "buyer1trace19"
Set-RedLabel "A92" "buyer2trace19"
Set-RedLabel "A93" "buyer3trace19"
Set-RedLabel "A94" "buyertracemx19"
Set-RedLabel "A95" "buyer1know19"
Set-RedLabel "A96" "buyer2know19"
Set-RedLabel "A97" "buyer3know19"
Set-RedLabel "A98" "buyerknowmx19"

# ---------------------------------------------------------------------
# Step 3: update view state (selection) to match the final saved state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B81").Select()
